$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-rank "Arabia Saudita" above "Filipinas" and "Mexico" (rows 37-39)
#    Row 37 becomes Arabia Saudita with freshly updated totals.
#    Row 38 becomes Filipinas, keeping its previous (unchanged) totals.
#    Row 39 becomes Mexico, keeping its previous (unchanged) totals.
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = "Arabia Saudita"
$ws.Range("B37").Value = 5862
$ws.Range("C37").Value = 493
$ws.Range("D37").Value = 931
$ws.Range("E37").Value = 4852
$ws.Range("F37").Value = 59
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 79

$ws.Range("A38").Value = "Filipinas"
$ws.Range("B38").Value = 5453
$ws.Range("C38").Value = 230
$ws.Range("D38").Value = 353
$ws.Range("E38").Value = 4751
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 14
$ws.Range("H38").Value = 349

$ws.Range("A39").Value = "Mexico"
$ws.Range("B39").Value = 5399
$ws.Range("C39").Value = 385
$ws.Range("D39").Value = 2125
$ws.Range("E39").Value = 2868
$ws.Range("F39").Value = 207
$ws.Range("G39").Value = 74
$ws.Range("H39").Value = 406

# ---------------------------------------------------------------------------
# 2) Re-rank "Barein" above "Hungria" (rows 63-64)
#    Row 63 becomes Barein with freshly updated totals.
#    Row 64 becomes Hungria, keeping its previous (unchanged) totals.
# ---------------------------------------------------------------------------
$ws.Range("A63").Value = "Barein"
$ws.Range("B63").Value = 1671
$ws.Range("C63").Value = 143
$ws.Range("D63").Value = 663
$ws.Range("E63").Value = 1001
$ws.Range("F63").Value = 3
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 7

$ws.Range("A64").Value = "Hungria"
$ws.Range("B64").Value = 1579
$ws.Range("C64").Value = 67
$ws.Range("D64").Value = 192
$ws.Range("E64").Value = 1253
$ws.Range("F64").Value = 58
$ws.Range("G64").Value = 12
$ws.Range("H64").Value = 134

# ---------------------------------------------------------------------------
# 3) Plain data refresh for Turquia (row 16) - no re-ranking needed
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = 26336
$ws.Range("C16").Value = 400
$ws.Range("E16").Value = 10422
$ws.Range("G16").Value = 40
$ws.Range("H16").Value = 1214

# ---------------------------------------------------------------------------
# 4) Plain data refresh for Argentina (row 54) - no re-ranking needed
# ---------------------------------------------------------------------------
$ws.Range("E54").Value = 1776
$ws.Range("G54").Value = 6
$ws.Range("H54").Value = 108

# ---------------------------------------------------------------------------
# 5) Plain data refresh for Uganda (row 150) - no re-ranking needed
# ---------------------------------------------------------------------------
$ws.Range("D150").Value = 12
$ws.Range("E150").Value = 43

# ---------------------------------------------------------------------------
# 6) Update the "last updated" timestamp banner
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 14:52"
